$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (A value, E value) reflecting the reordering of the
# "rank" (A) and "rate_value" (E) columns: within each 5-row block (same bank)
# rows are re-sorted by ascending rank.
$data = @{
    2 = @(24, 191.41)
    3 = @(25, 191.41)
    4 = @(26, 185.08)
    5 = @(27, 185.08)
    6 = @(28, 184.57)
    7 = @(29, 184.23)
    8 = @(30, 184.23)
    9 = @(31, 184.06)
    10 = @(32, 184.04)
    11 = @(33, 183.96)
    12 = @(7, 221.2)
    13 = @(8, 221.2)
    14 = @(9, 220.25)
    15 = @(10, 220.15)
    16 = @(11, 219.3)
    17 = @(14, 215.04)
    18 = @(15, 214.34)
    19 = @(16, 213.09)
    20 = @(17, 212.88)
    21 = @(18, 212.88)
    22 = @(34, 181.44)
    23 = @(35, 179.75)
    24 = @(36, 169.9)
    25 = @(37, 167.69)
    26 = @(38, 165.4)
    27 = @(5, 229.88)
    28 = @(6, 222.01)
    29 = @(19, 206.82)
    30 = @(21, 205.55)
    31 = @(22, 203.69)
    32 = @(0, 233.65)
    33 = @(1, 232.32)
    34 = @(2, 232.3)
    35 = @(3, 232.09)
    36 = @(4, 231.28)
    37 = @(12, 217.99)
    38 = @(13, 216.81)
    39 = @(20, 205.58)
    40 = @(23, 192.29)
    41 = @(69, 146.52)
    42 = @(63, 149.92)
    43 = @(64, 149.83)
    44 = @(65, 149.8)
    45 = @(66, 149.48)
    46 = @(67, 149.36)
    47 = @(58, 150.99)
    48 = @(59, 150.95)
    49 = @(60, 150.87)
    50 = @(61, 150.29)
    51 = @(62, 150.1)
    52 = @(41, 153.05)
    53 = @(44, 152.75)
    54 = @(48, 152.31)
    55 = @(56, 151.37)
    56 = @(57, 151.29)
    57 = @(39, 153.06)
    58 = @(40, 153.06)
    59 = @(42, 153.04)
    60 = @(43, 153.03)
    61 = @(45, 152.71)
    62 = @(51, 152.19)
    63 = @(52, 152.16)
    64 = @(53, 152)
    65 = @(54, 151.99)
    66 = @(55, 151.98)
    67 = @(46, 152.56)
    68 = @(47, 152.32)
    69 = @(49, 152.22)
    70 = @(50, 152.21)
    71 = @(68, 149.35)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 5).Value = $values[1]
}
